# Insert a new weekly data row at row 14 (pushing all subsequent rows down by one),
# containing the new price observation for Perejil (Terminal La Palmera de La Serena).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 14..156 down by inserting a brand-new row at position 14.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value = 8
$ws.Cells.Item(14, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(14, 3).Value = "Coquimbo"
$ws.Cells.Item(14, 4).Value = 44761
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = 100112044
$ws.Cells.Item(14, 7).Value = "Perejil"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 2400
$ws.Cells.Item(14, 11).Value = 2000
$ws.Cells.Item(14, 12).Value = 2500
$ws.Cells.Item(14, 13).Value = 2250
$ws.Cells.Item(14, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(14, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(14, 16).Value = 1500
$ws.Cells.Item(14, 17).Value = 1.5
$ws.Cells.Item(14, 18).Value = "Hortaliza"
